$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Workbook-level settings
# ---------------------------------------------------------------------------
# Remove workbook protection (no password was ever set; this just drops the
# now-unused <workbookProtection/> element from workbook.xml)
$wb.Unprotect()

# ---------------------------------------------------------------------------
# Frutas sheet: fix the "for loop" off-by-one that hard coded "Uva" as the
# 5th fruit instead of using the loop placeholder, and dress up the table
# with centered text + a grid of thin borders + a bold header row.
# ---------------------------------------------------------------------------
$frutas = $wb.Worksheets.Item("Frutas")

$frutas.Range("A5").Value2 = "Fruta 1"

$frutas.Columns.Item(1).ColumnWidth = 8.43
$frutas.Columns.Item(2).ColumnWidth = 11.0
$frutas.Columns.Item(3).ColumnWidth = 11.57
$frutas.Columns.Item(4).ColumnWidth = 8.43

$frutasAll = $frutas.Range("A1:C5")
$frutasAll.HorizontalAlignment = -4108

$frutasHeader = $frutas.Range("A1:C1")
$frutasHeader.Borders.LineStyle = 1
$frutasHeader.Font.Bold = $true

$frutasBody = $frutas.Range("A2:C5")
$frutasBody.Borders.LineStyle = 1

$null = $frutas.Range("F7").Select()

# ---------------------------------------------------------------------------
# Cores sheet: brand-new sheet, placed right after Frutas, mirroring the
# same header/body styling used on the Frutas table.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cores = $wb.Worksheets.Add($null, $lastSheet)
$cores.Name = "Cores"

$cores.Range("A1").Value2 = "Primária"
$cores.Range("B1").Value2 = "Secundária"
$cores.Range("A2").Value2 = "Azul"
$cores.Range("B2").Value2 = "Amarelo"
$cores.Range("A3").Value2 = "Verde"
$cores.Range("B3").Value2 = "Vermelho"
$cores.Range("A4").Value2 = "Rosa"
$cores.Range("B4").Value2 = "Preto"
$cores.Range("A5").Value2 = "Branco"
$cores.Range("B5").Value2 = "Laranja"

$cores.Columns.Item(1).ColumnWidth = 8.43
$cores.Columns.Item(2).ColumnWidth = 10.14
$cores.Columns.Item(3).ColumnWidth = 8.43

$coresAll = $cores.Range("A1:B5")
$coresAll.HorizontalAlignment = -4108

$coresHeader = $cores.Range("A1:B1")
$coresHeader.Borders.LineStyle = 1
$coresHeader.Font.Bold = $true

$coresBody = $cores.Range("A2:B5")
$coresBody.Borders.LineStyle = 1

$null = $cores.Range("A1:B1").Select()

# ---------------------------------------------------------------------------
# Make Frutas the active tab again (adding Cores activated it by default)
# ---------------------------------------------------------------------------
$frutas.Activate()

Write-Output "done"
